$d = $word.ActiveDocument

# The diff removes the explicit bold run formatting (<w:rPr><w:bCs/><w:b/></w:rPr>)
# from the two header-row cells ("Content" / "Web Site") of the table, while
# leaving the paragraph formatting (style/alignment) and the text untouched.
#
# This runtime's Font.Bold / Range.Bold setters always materialize an explicit
# "off" value (<w:b w:val="0"/>) instead of deleting the <w:b/>/<w:bCs/> nodes,
# so we rebuild the run via InsertXML (which replaces a Range's contents with
# exactly the OOXML we supply) using a bare <w:r> with no <w:rPr> at all.

function Clear-RunBoldFormatting($cell) {
    $cellRange = $cell.Range
    $para = $cellRange.Paragraphs.Item(1)

    # Preserve the paragraph's existing style/alignment.
    $styleName = $para.Range.ParagraphStyle.NameLocal
    $alignment = $para.Alignment

    switch ($alignment) {
        1 { $jc = "center" }
        2 { $jc = "right" }
        3 { $jc = "both" }
        default { $jc = "left" }
    }

    $text = $cellRange.Text
    # A cell's Range.Text ends with the paragraph mark (CR, chr 13) followed
    # by the end-of-cell mark (chr 7); drop both so we only touch the actual
    # run text.
    $text = $text.Substring(0, $text.Length - 2)

    $target = $d.Range($cellRange.Start, $cellRange.Start + $text.Length)

    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:pPr><w:pStyle w:val="' + $styleName + '"/><w:jc w:val="' + $jc + '"/></w:pPr>' +
           '<w:r><w:t xml:space="preserve">' + $text + '</w:t></w:r></w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}

$table = $d.Tables.Item(1)
$headerRow = $table.Rows.Item(1)

Clear-RunBoldFormatting $headerRow.Cells.Item(1)
Clear-RunBoldFormatting $headerRow.Cells.Item(2)
